$wb = $excel.ActiveWorkbook

# Add the new "AffiliatedCompanies" worksheet right after "TopRelationships"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "AffiliatedCompanies"

# Populate the field labels for the new sheet
$newSheet.Range("A1").Value = "Company Name:"
$newSheet.Range("A2").Value = "Type:"
$newSheet.Range("A3").Value = "Company Type:"

# Make the new sheet the active tab with the same selected cell as the source edit
$newSheet.Range("K14").Select()
